# chore: update Sheets via scheduled runner
# Refreshes marketboard price / profit figures (columns H-N) across the
# ALC, ARM, BSM, CRP, CUL, LTW and WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M7").ClearContents()
$ws.Range("H7").Value = 2750
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 2750
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 2750
$ws.Range("N7").Value = -2974

$ws.Range("M14").ClearContents()
$ws.Range("H14").Value = 2750
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 2750
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 2750
$ws.Range("N14").Value = -3132

$ws.Range("H64").Value = 3358.3333
$ws.Range("I64").Value = 3240
$ws.Range("J64").Value = 3385.2273
$ws.Range("K64").Value = 3240
$ws.Range("L64").Value = 3385.2273
$ws.Range("M64").Value = -2992
$ws.Range("N64").Value = -3881.2273

$ws.Range("H67").Value = 3358.3333
$ws.Range("I67").Value = 3240
$ws.Range("J67").Value = 3385.2273
$ws.Range("K67").Value = 3240
$ws.Range("L67").Value = 3385.2273
$ws.Range("M67").Value = -2382
$ws.Range("N67").Value = -5101.2273

$ws.Range("H69").Value = 3816.6667
$ws.Range("J69").Value = 3816.6667
$ws.Range("L69").Value = 11450.0001
$ws.Range("N69").Value = -13198.0001

$ws.Range("H70").Value = 2078.1738
$ws.Range("I70").Value = 2028.0312
$ws.Range("J70").Value = 2192.7856
$ws.Range("K70").Value = 6084.0936
$ws.Range("L70").Value = 6578.3568
$ws.Range("M70").Value = -5814.0936
$ws.Range("N70").Value = -7118.3568

$ws.Range("H72").Value = 3816.6667
$ws.Range("J72").Value = 3816.6667
$ws.Range("L72").Value = 34350.0003
$ws.Range("N72").Value = -43086.0003

$ws.Range("H73").Value = 2078.1738
$ws.Range("I73").Value = 2028.0312
$ws.Range("J73").Value = 2192.7856
$ws.Range("K73").Value = 6084.0936
$ws.Range("L73").Value = 6578.3568
$ws.Range("M73").Value = -5148.0936
$ws.Range("N73").Value = -8450.356800000001

$ws.Range("H96").Value = 340.25
$ws.Range("I96").Value = 285.0476
$ws.Range("J96").Value = 445.63635
$ws.Range("K96").Value = 855.1428
$ws.Range("L96").Value = 1336.90905
$ws.Range("M96").Value = 517.8572
$ws.Range("N96").Value = -4082.90905

$ws.Range("H98").Value = 1282.963
$ws.Range("I98").Value = 864.1667
$ws.Range("J98").Value = 4633.3335
$ws.Range("K98").Value = 864.1667
$ws.Range("L98").Value = 4633.3335
$ws.Range("M98").Value = 633.8333
$ws.Range("N98").Value = -7629.3335

$ws.Range("H100").Value = 33334988
$ws.Range("I100").Value = 1618.091
$ws.Range("K100").Value = 1618.091
$ws.Range("M100").Value = -1077.091

$ws.Range("H103").Value = 9524435
$ws.Range("J103").Value = 15385319
$ws.Range("L103").Value = 46155957
$ws.Range("N103").Value = -46157129

$ws.Range("H122").Value = 1282.963
$ws.Range("I122").Value = 864.1667
$ws.Range("J122").Value = 4633.3335
$ws.Range("K122").Value = 2592.5001
$ws.Range("L122").Value = 13900.0005
$ws.Range("M122").Value = -142.5001000000002
$ws.Range("N122").Value = -18800.0005

$ws.Range("H128").Value = 36160
$ws.Range("J128").Value = 36160
$ws.Range("L128").Value = 36160
$ws.Range("N128").Value = -46120

$ws.Range("H133").Value = 39606.168
$ws.Range("J133").Value = 39606.168
$ws.Range("L133").Value = 39606.168
$ws.Range("N133").Value = -49726.168

$ws.Range("H134").Value = 47743.7
$ws.Range("J134").Value = 47743.7
$ws.Range("L134").Value = 47743.7
$ws.Range("N134").Value = -57883.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32016
$ws.Range("I32").Value = 9247.083000000001
$ws.Range("K32").Value = 9247.083000000001
$ws.Range("M32").Value = -8960.083000000001

$ws.Range("H36").Value = 250001340
$ws.Range("I36").Value = 1785.3334
$ws.Range("J36").Value = 1000000000
$ws.Range("K36").Value = 1785.3334
$ws.Range("L36").Value = 1000000000
$ws.Range("M36").Value = -1439.3334
$ws.Range("N36").Value = -1000000692

$ws.Range("H39").Value = 15703.667
$ws.Range("I39").Value = 9772
$ws.Range("K39").Value = 9772
$ws.Range("M39").Value = -9252

$ws.Range("H103").Value = 43944
$ws.Range("J103").Value = 43944
$ws.Range("L103").Value = 43944
$ws.Range("N103").Value = -46288

$ws.Range("H109").Value = 19900
$ws.Range("J109").Value = 19900
$ws.Range("L109").Value = 19900
$ws.Range("N109").Value = -22674

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 25298
$ws.Range("J19").Value = 25298
$ws.Range("L19").Value = 25298
$ws.Range("N19").Value = -25644

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 48307.8
$ws.Range("J20").Value = 48307.8
$ws.Range("L20").Value = 48307.8
$ws.Range("N20").Value = -48779.8

$ws.Range("H30").Value = 48307.8
$ws.Range("J30").Value = 48307.8
$ws.Range("L30").Value = 48307.8
$ws.Range("N30").Value = -48489.8

$ws.Range("H128").Value = 48307.8
$ws.Range("J128").Value = 48307.8
$ws.Range("L128").Value = 48307.8
$ws.Range("N128").Value = -58267.8

$ws.Range("H135").Value = 40416.19
$ws.Range("J135").Value = 40416.19
$ws.Range("L135").Value = 40416.19
$ws.Range("N135").Value = -50556.19

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2080.7144
$ws.Range("I137").Value = 1856.1538
$ws.Range("K137").Value = 5568.4614
$ws.Range("M137").Value = -468.4614000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1531.579
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1776.9231
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1776.9231
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -2366.9231

$ws.Range("H27").Value = 1531.579
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1776.9231
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1776.9231
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1990.9231

$ws.Range("H35").Value = 1098.4
$ws.Range("I35").Value = 1098.4
$ws.Range("K35").Value = 1098.4
$ws.Range("M35").Value = -762.4000000000001

$ws.Range("H63").Value = 26396.25
$ws.Range("J63").Value = 26396.25
$ws.Range("L63").Value = 26396.25
$ws.Range("N63").Value = -27894.25

$ws.Range("H66").Value = 26396.25
$ws.Range("J66").Value = 26396.25
$ws.Range("L66").Value = 79188.75
$ws.Range("N66").Value = -86676.75

$ws.Range("H108").Value = 15067.667
$ws.Range("J108").Value = 15067.667
$ws.Range("L108").Value = 15067.667
$ws.Range("N108").Value = -22747.667

$ws.Range("H127").Value = 49966.8
$ws.Range("J127").Value = 49966.8
$ws.Range("L127").Value = 49966.8
$ws.Range("N127").Value = -59886.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 527436.5
$ws.Range("I100").Value = 1267.1428
$ws.Range("J100").Value = 2000710.6
$ws.Range("K100").Value = 2534.2856
$ws.Range("L100").Value = 4001421.2
$ws.Range("M100").Value = -1993.2856
$ws.Range("N100").Value = -4002503.2

$ws.Range("H107").Value = 3105
$ws.Range("I107").Value = 1681.4546
$ws.Range("J107").Value = 5062.375
$ws.Range("K107").Value = 5044.3638
$ws.Range("L107").Value = 15187.125
$ws.Range("M107").Value = -3124.3638
$ws.Range("N107").Value = -19027.125

$ws.Range("H108").Value = 25123.334
$ws.Range("J108").Value = 25123.334
$ws.Range("L108").Value = 25123.334
$ws.Range("N108").Value = -32803.334

$ws.Range("H113").Value = 53686.79
$ws.Range("I113").Value = 71807.07000000001
$ws.Range("J113").Value = 2950
$ws.Range("K113").Value = 215421.21
$ws.Range("L113").Value = 8850
$ws.Range("M113").Value = -213251.21
$ws.Range("N113").Value = -13190
